$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "Architecture" (old E) for "Inapplicable CPE URIs".
# This shifts old E..K (Architecture..Remarks) to F..L.
$ws.Columns.Item(5).Insert()

# --- Header row (row 1) ---
$ws.Range("E1").Value = "Inapplicable CPE URIs"

# --- Row 2: Windows 11 25H2 ---
$ws.Range("D2").Value = "cpe:/o:microsoft:windows_11_25h2"
$ws.Range("E2").Value = "cpe:/o:microsoft:windows_11"
$ws.Range("L2").Value = "cpe:/o:microsoft:windows_11 is considered a outdated for this version of windows."

# --- Row 3: Windows 11 24H2 ---
$ws.Range("D3").Value = "cpe:/o:microsoft:windows_11_24h2"
$ws.Range("E3").Value = "cpe:/o:microsoft:windows_11"
$ws.Range("L3").Value = "cpe:/o:microsoft:windows_11 is considered a outdated for this version of windows."

# --- Row 4: Windows 11 23H2 ---
$ws.Range("D4").Value = "cpe:/o:microsoft:windows_11_23h2"
$ws.Range("E4").Value = "cpe:/o:microsoft:windows_11"
$ws.Range("L4").Value = "cpe:/o:microsoft:windows_11 is considered a outdated for this version of windows."

# --- Column widths (approximate best-fit look) ---
$ws.Columns.Item(1).ColumnWidth = 13.998697916666666
$ws.Columns.Item(2).ColumnWidth = 11.830729166666666
$ws.Columns.Item(3).ColumnWidth = 13.830729166666666
$ws.Columns.Item(4).ColumnWidth = 27.666666666666668
$ws.Columns.Item(5).ColumnWidth = 22.666666666666668
$ws.Columns.Item(6).ColumnWidth = 12.330729166666666
$ws.Columns.Item(7).ColumnWidth = 13.666666666666666
$ws.Columns.Item(9).ColumnWidth = 12.998697916666666
$ws.Columns.Item(10).ColumnWidth = 7.498697916666667
$ws.Columns.Item(11).ColumnWidth = 17.498697916666668
$ws.Columns.Item(12).ColumnWidth = 62.330729166666664

# --- AutoFilter: extend from A1:K65001 to A1:L65001 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:L65001").AutoFilter() | Out-Null

# --- Defined name _FilterDatabase: extend to column L ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Artifacts!_FilterDatabase") {
        $n.RefersTo = '=Artifacts!$A$1:$L$65001'
    }
}

# --- Selection ---
$ws.Range("K7").Select() | Out-Null

Write-Output "done"
